$wb = $excel.ActiveWorkbook

# Data refresh: Sheets/Typhon_Profits.xlsx (scheduled runner update)
# Each entry: sheet name, row, and the column letter -> new value map.
# A $null value means the cell is cleared (removed) entirely.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2171.5833
$ws.Range("I51").Value = 1819.8
$ws.Range("K51").Value = 1819.8
$ws.Range("M51").Value = -1335.8

$ws.Range("H80").Value = 8703125
$ws.Range("I80").Value = 629.7143
$ws.Range("J80").Value = 17405620
$ws.Range("K80").Value = 1889.1429
$ws.Range("L80").Value = 52216860
$ws.Range("M80").Value = -891.1428999999998
$ws.Range("N80").Value = -52218856

$ws.Range("H83").Value = 8703125
$ws.Range("I83").Value = 629.7143
$ws.Range("J83").Value = 17405620
$ws.Range("K83").Value = 5667.428699999999
$ws.Range("L83").Value = 156650580
$ws.Range("M83").Value = -675.4286999999995
$ws.Range("N83").Value = -156660564

$ws.Range("H129").Value = 1154.6
$ws.Range("I129").Value = 478.57144
$ws.Range("K129").Value = 1435.71432
$ws.Range("M129").Value = 3564.28568

$ws.Range("H132").Value = 1803.3684
$ws.Range("I132").Value = 1864.6666
$ws.Range("K132").Value = 5593.9998
$ws.Range("M132").Value = -3063.9998

$ws.Range("H137").Value = 1173.8772
$ws.Range("I137").Value = 1071.3654
$ws.Range("K137").Value = 3214.0962
$ws.Range("M137").Value = -664.0962

$ws.Range("H138").Value = 2599.4827
$ws.Range("I138").Value = 4840
$ws.Range("J138").Value = 2388.1133
$ws.Range("K138").Value = 14520
$ws.Range("L138").Value = 7164.3399
$ws.Range("M138").Value = -9380
$ws.Range("N138").Value = -17444.3399

$ws.Range("H141").Value = 1316.4615
$ws.Range("I141").Value = 1096.9131
$ws.Range("J141").Value = 2999.6667
$ws.Range("K141").Value = 3290.7393
$ws.Range("L141").Value = 8999.000100000001
$ws.Range("M141").Value = 1889.2607
$ws.Range("N141").Value = -19359.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4088.3635
$ws.Range("I32").Value = 4707.7554
$ws.Range("J32").Value = 1301.1
$ws.Range("K32").Value = 4707.7554
$ws.Range("L32").Value = 1301.1
$ws.Range("M32").Value = -4420.7554
$ws.Range("N32").Value = -1875.1

$ws.Range("H132").Value = 13839.609
$ws.Range("I132").Value = 1412.1714
$ws.Range("K132").Value = 4236.5142
$ws.Range("M132").Value = -1706.5142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1598
$ws.Range("I86").Value = 1417.3334
$ws.Range("J86").Value = 2085.8
$ws.Range("K86").Value = 1417.3334
$ws.Range("L86").Value = 2085.8
$ws.Range("M86").Value = -294.3334
$ws.Range("N86").Value = -4331.8

$ws.Range("H89").Value = 1598
$ws.Range("I89").Value = 1417.3334
$ws.Range("J89").Value = 2085.8
$ws.Range("K89").Value = 7086.666999999999
$ws.Range("L89").Value = 10429
$ws.Range("M89").Value = -1470.666999999999
$ws.Range("N89").Value = -21661

$ws.Range("H134").Value = 3237.738
$ws.Range("I134").Value = 3579.3235
$ws.Range("J134").Value = 1786
$ws.Range("K134").Value = 10737.9705
$ws.Range("L134").Value = 5358
$ws.Range("M134").Value = -8202.970499999999
$ws.Range("N134").Value = -10428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2378.5107
$ws.Range("I132").Value = 1798.4706
$ws.Range("J132").Value = 3895.5386
$ws.Range("K132").Value = 5395.4118
$ws.Range("L132").Value = 11686.6158
$ws.Range("M132").Value = -2865.4118
$ws.Range("N132").Value = -16746.6158

$ws.Range("H134").Value = 884.2414
$ws.Range("I134").Value = 767.5217
$ws.Range("J134").Value = 1331.6666
$ws.Range("K134").Value = 2302.5651
$ws.Range("L134").Value = 3994.9998
$ws.Range("M134").Value = 232.4349000000002
$ws.Range("N134").Value = -9064.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 4841.6665
$ws.Range("J76").Value = 4841.6665
$ws.Range("L76").Value = 14524.9995
$ws.Range("N76").Value = -15290.9995

$ws.Range("H79").Value = 4841.6665
$ws.Range("J79").Value = 4841.6665
$ws.Range("L79").Value = 14524.9995
$ws.Range("N79").Value = -17176.9995

$ws.Range("H129").Value = 295000.53
$ws.Range("I129").Value = 776.6667
$ws.Range("J129").Value = 358048.5
$ws.Range("K129").Value = 2330.0001
$ws.Range("L129").Value = 1074145.5
$ws.Range("M129").Value = 2669.9999
$ws.Range("N129").Value = -1084145.5

$ws.Range("H131").Value = 793.37
$ws.Range("I131").Value = 687.8
$ws.Range("J131").Value = 798.92633
$ws.Range("K131").Value = 2063.4
$ws.Range("L131").Value = 2396.77899
$ws.Range("M131").Value = 2976.6
$ws.Range("N131").Value = -12476.77899

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 103500
$ws.Range("J20").Value = 103500
$ws.Range("L20").Value = 103500
$ws.Range("N20").Value = -103990

$ws.Range("I24").Value = 200000
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 200000
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -199827
$ws.Range("N24").ClearContents()

$ws.Range("H31").Value = 4031
$ws.Range("I31").Value = 4031
$ws.Range("K31").Value = 4031
$ws.Range("M31").Value = -3739

$ws.Range("H37").Value = 4031
$ws.Range("I37").Value = 4031
$ws.Range("K37").Value = 4031
$ws.Range("M37").Value = -3754

$ws.Range("H123").Value = 10253.333
$ws.Range("J123").Value = 10253.333
$ws.Range("L123").Value = 10253.333
$ws.Range("N123").Value = -15153.333

$ws.Range("H132").Value = 47816.082
$ws.Range("I132").Value = 8828
$ws.Range("J132").Value = 102399.4
$ws.Range("K132").Value = 26484
$ws.Range("L132").Value = 307198.2
$ws.Range("M132").Value = -23954
$ws.Range("N132").Value = -312258.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3436.7273
$ws.Range("I93").Value = 3412.5
$ws.Range("K93").Value = 3412.5
$ws.Range("M93").Value = -2164.5

$ws.Range("H139").Value = 59999.5
$ws.Range("J139").Value = 59999.5
$ws.Range("L139").Value = 59999.5
$ws.Range("N139").Value = -70279.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 69661.336
$ws.Range("I14").Value = 102002
$ws.Range("K14").Value = 102002
$ws.Range("M14").Value = -101834

$ws.Range("H30").Value = 1666.3334
$ws.Range("I30").Value = 1666.3334
$ws.Range("K30").Value = 1666.3334
$ws.Range("M30").Value = -1559.3334

$ws.Range("H81").Value = 111112670
$ws.Range("I81").Value = 1755.75
$ws.Range("J81").Value = 1000000000
$ws.Range("K81").Value = 3511.5
$ws.Range("L81").Value = 2000000000
$ws.Range("M81").Value = -2450.5
$ws.Range("N81").Value = -2000002122

$ws.Range("H84").Value = 111112670
$ws.Range("I84").Value = 1755.75
$ws.Range("J84").Value = 1000000000
$ws.Range("K84").Value = 17557.5
$ws.Range("L84").Value = 10000000000
$ws.Range("M84").Value = -12253.5
$ws.Range("N84").Value = -10000010608
